# Añadido story/tareas para borrar ingredientes del inventario
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark existing sections "CREACIÓN DE UNA CUENTA" (row 2) and
#     "AÑADIR RECETAS A LA BASE DE DATOS" (row 36) as HECHO, matching the
#     formatting already used for the other completed sections (e.g. C15).
$ws.Range("C15").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C36").PasteSpecial(-4122)
$ws.Range("C15").Copy()
$excel.CutCopyMode = $false

# --- New story: "BORRAR INGREDIENTES" -----------------------------------
# Section header row (copy formatting from the "AÑADIR RECETAS A LA BASE
# DE DATOS" header row, which also has no HECHO flag).
$ws.Range("A36:B36").Copy()
$ws.Range("A46:B46").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A46").Value = 2
$ws.Range("B46").Value = "BORRAR INGREDIENTES"

# Sub-task rows: copy formatting from row 45 (a similarly-styled task row)
$ws.Range("A45:F45").Copy()
$ws.Range("A47:F47").PasteSpecial(-4122)
$ws.Range("A48:F48").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B47").Value = "Crear una opción en el menú para que el usuario borre ingredientes de su inventario"
$ws.Range("D47").Value = "Enrique"
$ws.Range("E47").Value = 0.5
$ws.Range("F47").ClearContents()

$ws.Range("B48").Value = "Crear la consulta sql para borrar un ingrediente introducido por el usuario"
$ws.Range("D48").Value = "Enrique"
$ws.Range("E48").Value = 0.1
$ws.Range("F48").ClearContents()

# --- View state: zoom + selection ---------------------------------------
$excel.ActiveWindow.Zoom = 125
$ws.Range("D54").Select()
